# Apply updated cryptos list values (Price column D, Volume(1h) column E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.526.06"
$ws.Range("E2").Value = "  +5.28%  "

$ws.Range("D3").Value = "1.724.60"
$ws.Range("E3").Value = "  +4.08%  "

$ws.Range("D4").Value = "'" + "1.004"
$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "'" + "225.87"
$ws.Range("E5").Value = "  +3.26%  "

$ws.Range("D6").Value = "'" + "0.5371"
$ws.Range("E6").Value = "  +2.46%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").Value = "'" + "0.2676"
$ws.Range("E8").Value = "  +0.68%  "

$ws.Range("D9").Value = "'" + "0.06601"
$ws.Range("E9").Value = "  +3.81%  "

$ws.Range("D10").Value = "'" + "21.74"
$ws.Range("E10").Value = "  +5.61%  "

$ws.Range("D11").Value = "'" + "0.07744"
$ws.Range("E11").Value = "  +0.66%  "

$ws.Range("D12").Value = "'" + "4.619"
$ws.Range("E12").Value = "  +0.19%  "

$ws.Range("D13").Value = "1.720.59"
$ws.Range("E13").Value = "  +2.10%  "

$ws.Range("D14").Value = "1.961.69"
$ws.Range("E14").Value = "  +4.08%  "

$ws.Range("D15").Value = "'" + "0.5867"
$ws.Range("E15").Value = "  +4.21%  "

$ws.Range("D16").Value = "0.0₅8309"
$ws.Range("E16").Value = "  +1.34%  "

$ws.Range("D17").Value = "'" + "68.00"
$ws.Range("E17").Value = "  +3.81%  "

$ws.Range("D18").Value = "27.547.48"
$ws.Range("E18").Value = "  +5.36%  "

$ws.Range("D19").Value = "'" + "221.76"
$ws.Range("E19").Value = "  +15.15%  "

$ws.Range("D20").Value = "'" + "1.004"
$ws.Range("E20").Value = "  +0.08%  "

$ws.Range("D21").Value = "'" + "4.743"
$ws.Range("E21").Value = "  +1.72%  "

$ws.Range("E22").Value = "  +1.54%  "

$ws.Range("D23").Value = "'" + "6.095"
$ws.Range("E23").Value = "  +2.18%  "

$ws.Range("E24").Value = "  +0.09%  "

$ws.Range("D25").Value = "'" + "148.15"
$ws.Range("E25").Value = "  +1.97%  "

$ws.Range("D26").Value = "'" + "1.695"
$ws.Range("E26").Value = "  +12.09%  "

$ws.Range("D27").Value = "'" + "0.1233"
$ws.Range("E27").Value = "  +2.82%  "

$ws.Range("D28").Value = "'" + "7.398"
$ws.Range("E28").Value = "  +1.76%  "

$ws.Range("D29").Value = "'" + "16.67"
$ws.Range("E29").Value = "  +4.34%  "

$ws.Range("D30").Value = "'" + "0.05541"
$ws.Range("E30").Value = "  +1.32%  "

$ws.Range("E31").Value = "  +2.41%  "

$ws.Range("D32").Value = "'" + "3.544"
$ws.Range("E32").Value = "  +2.17%  "

$ws.Range("D33").Value = "'" + "3.461"
$ws.Range("E33").Value = "  +2.56%  "

$ws.Range("D34").Value = "'" + "1.662"
$ws.Range("E34").Value = "  +6.08%  "

$ws.Range("D35").Value = "'" + "0.9602"
$ws.Range("E35").Value = "  +0.54%  "

$ws.Range("E36").Value = "  +1.52%  "

$ws.Range("E37").Value = "  +1.80%  "

$ws.Range("D38").Value = "'" + "0.5948"
$ws.Range("E38").Value = "  +4.57%  "

$ws.Range("D39").Value = "'" + "0.01648"
$ws.Range("E39").Value = "  +3.77%  "

$ws.Range("D40").Value = "'" + "5.933"
$ws.Range("E40").Value = "  +0.87%  "

$ws.Range("D41").Value = "1.058.93"
$ws.Range("E41").Value = "  +2.93%  "

$ws.Range("D42").Value = "'" + "0.8547"
$ws.Range("E42").Value = "  +2.57%  "

$ws.Range("D44").Value = "'" + "101.56"
$ws.Range("E44").Value = "  +0.19%  "

$ws.Range("D45").Value = "1.867.92"
$ws.Range("E45").Value = "  +4.01%  "

$ws.Range("E46").Value = "  +5.56%  "

$ws.Range("D47").Value = "'" + "59.05"
$ws.Range("E47").Value = "  +2.09%  "

$ws.Range("D48").Value = "'" + "8.191"
$ws.Range("E48").Value = "  +1.85%  "

$ws.Range("E49").Value = "  +2.25%  "

$ws.Range("D50").Value = "'" + "1.001"
$ws.Range("E50").Value = "  -0.06%  "

$ws.Range("D51").Value = "'" + "0.05276"
$ws.Range("E51").Value = "  +1.67%  "
